$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.17532746519896136,
    -0.0099999995554611587,
    -0.0089999995563676549,
    0.20999999883007092,
    -0.0059999995664359318,
    -0.0059999995485391366,
    -0.019999999468049978,
    -0.048139182020049631,
    -0.0059999995407382656,
    -0.0059999995396182726,
    -0.0044999995479706456,
    -0.0059999995395481065,
    -0.0059999995403128281,
    -0.011999999507151138,
    0.045457563862300709,
    -0.0059999995408310802,
    -0.0059999995389210525,
    -0.0089999995215999107,
    -0.047188376518308761,
    -0.0089999995562397572,
    -0.0089999995556766521,
    -0.0089999995552707546,
    -0.0089999995525884557,
    -0.041999999359275009,
    -0.041999999355874174,
    -0.0059999995465958023,
    -0.005999999544219925,
    -0.00599999953410979,
    -0.011999999493504276,
    -0.019999999445131422,
    -0.014999999468427561,
    -0.0209999994337986,
    -0.0059999995180826105
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
